$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 139, shifting existing rows 139-179 down to 140-180.
$ws.Rows.Item(139).Insert()

# Populate the newly inserted row 139 with the new weekly entry.
$ws.Range("A139").Value2 = 4
$ws.Range("B139").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C139").Value2 = "Los Lagos"
$ws.Range("D139").Value2 = 44642
$ws.Range("E139").Value2 = 10
$ws.Range("F139").Value2 = "Fruta"
$ws.Range("G139").Value2 = 100108
$ws.Range("H139").Value2 = "Tropicales y subtropicales"
$ws.Range("I139").Value2 = 100108002
$ws.Range("J139").Value2 = "Mango"
$ws.Range("K139").Value2 = "Sin especificar"
$ws.Range("L139").Value2 = "Primera"
$ws.Range("M139").Value2 = 200
$ws.Range("N139").Value2 = 8000
$ws.Range("O139").Value2 = 8500
$ws.Range("P139").Value2 = 8250
$ws.Range("Q139").Value2 = "$/bandeja 4 kilos"
$ws.Range("R139").Value2 = "Perú"
$ws.Range("S139").Value2 = 2062
$ws.Range("T139").Value2 = 4
